$wb = $excel.ActiveWorkbook

# Sheet: ALC (19 cells)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 242.28572
$ws.Range("I6").Value = 249.33333
$ws.Range("K6").Value = 747.99999
$ws.Range("M6").Value = -635.99999
$ws.Range("H107").Value = 38463856
$ws.Range("I107").Value = 45455660
$ws.Range("J107").Value = 8938
$ws.Range("K107").Value = 45455660
$ws.Range("L107").Value = 8938
$ws.Range("M107").Value = -45453740
$ws.Range("N107").Value = -12778
$ws.Range("H137").Value = 5324456
$ws.Range("I137").Value = 11907810
$ws.Range("K137").Value = 35723430
$ws.Range("M137").Value = -35720880
$ws.Range("H138").Value = 3707.6453
$ws.Range("J138").Value = 4572.5625
$ws.Range("L138").Value = 13717.6875
$ws.Range("N138").Value = -23997.6875

# Sheet: ARM (24 cells)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 11503
$ws.Range("J3").Value = 11503
$ws.Range("L3").Value = 11503
$ws.Range("N3").Value = -11733
$ws.Range("H4").Value = 41.666668
$ws.Range("H5").Value = 355.73334
$ws.Range("J5").Value = 134.2
$ws.Range("L5").Value = 134.2
$ws.Range("N5").Value = -358.2
$ws.Range("H8").Value = 1997.5
$ws.Range("J8").Value = 1997.5
$ws.Range("L8").Value = 1997.5
$ws.Range("N8").Value = -2285.5
$ws.Range("H102").Value = 1341.1333
$ws.Range("I102").Value = 1294.2142
$ws.Range("K102").Value = 1294.2142
$ws.Range("M102").Value = 327.7858000000001
$ws.Range("H122").Value = 2634.4
$ws.Range("I122").Value = 1884.6666
$ws.Range("J122").Value = 3134.2222
$ws.Range("K122").Value = 5653.9998
$ws.Range("L122").Value = 9402.6666
$ws.Range("M122").Value = -3203.9998
$ws.Range("N122").Value = -14302.6666

# Sheet: BSM (33 cells)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 355.73334
$ws.Range("J4").Value = 134.2
$ws.Range("L4").Value = 134.2
$ws.Range("N4").Value = -364.2
$ws.Range("H86").Value = 334451
$ws.Range("I86").Value = 1341.2
$ws.Range("J86").Value = 2000000
$ws.Range("K86").Value = 1341.2
$ws.Range("L86").Value = 2000000
$ws.Range("M86").Value = -218.2
$ws.Range("N86").Value = -2002246
$ws.Range("H89").Value = 334451
$ws.Range("I89").Value = 1341.2
$ws.Range("J89").Value = 2000000
$ws.Range("K89").Value = 6706
$ws.Range("L89").Value = 10000000
$ws.Range("M89").Value = -1090
$ws.Range("N89").Value = -10011232
$ws.Range("H94").Value = 1348.5454
$ws.Range("I94").Value = 1372
$ws.Range("K94").Value = 1372
$ws.Range("M94").Value = -921
$ws.Range("H105").Value = 34493770
$ws.Range("J105").Value = 2527.8
$ws.Range("L105").Value = 2527.8
$ws.Range("N105").Value = -6021.8
$ws.Range("H134").Value = 3218.4055
$ws.Range("I134").Value = 2117.5334
$ws.Range("J134").Value = 7936.4287
$ws.Range("K134").Value = 6352.600199999999
$ws.Range("L134").Value = 23809.2861
$ws.Range("M134").Value = -3817.600199999999
$ws.Range("N134").Value = -28879.2861

# Sheet: CRP (34 cells)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 456.5
$ws.Range("I7").Value = 459.11765
$ws.Range("K7").Value = 459.11765
$ws.Range("M7").Value = -346.11765
$ws.Range("H22").Value = 2635.75
$ws.Range("I22").Value = 771.5
$ws.Range("J22").Value = 4500
$ws.Range("K22").Value = 771.5
$ws.Range("L22").Value = 4500
$ws.Range("M22").Value = -421.5
$ws.Range("N22").Value = -5200
$ws.Range("H93").Value = 9608.909
$ws.Range("I93").Value = 2857.3333
$ws.Range("K93").Value = 2857.3333
$ws.Range("M93").Value = -985.3332999999998
$ws.Range("H103").Value = 34509.2
$ws.Range("J103").Value = 79996
$ws.Range("L103").Value = 79996
$ws.Range("N103").Value = -82340
$ws.Range("H104").Value = 37997.4
$ws.Range("J104").Value = 59995.668
$ws.Range("L104").Value = 59995.668
$ws.Range("N104").Value = -65237.668
$ws.Range("H105").Value = 3810.3125
$ws.Range("I105").Value = 5567.6665
$ws.Range("K105").Value = 5567.6665
$ws.Range("M105").Value = -3820.6665
$ws.Range("H132").Value = 4157
$ws.Range("I132").Value = 3477
$ws.Range("J132").Value = 9257
$ws.Range("K132").Value = 10431
$ws.Range("L132").Value = 27771
$ws.Range("M132").Value = -7901
$ws.Range("N132").Value = -32831

# Sheet: CUL (39 cells)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 638
$ws.Range("I5").Value = 531.7143
$ws.Range("J5").Value = 720.6667
$ws.Range("K5").Value = 1595.1429
$ws.Range("L5").Value = 2162.0001
$ws.Range("M5").Value = -1483.1429
$ws.Range("N5").Value = -2386.0001
$ws.Range("H7").Value = 48.8
$ws.Range("I7").Value = 49.666668
$ws.Range("J7").Value = 47.5
$ws.Range("K7").Value = 149.000004
$ws.Range("L7").Value = 142.5
$ws.Range("M7").Value = -37.00000399999999
$ws.Range("N7").Value = -366.5
$ws.Range("H14").Value = 1848.4667
$ws.Range("I14").Value = 1848.4667
$ws.Range("K14").Value = 5545.4001
$ws.Range("M14").Value = -5372.4001
$ws.Range("H92").Value = 1952.1538
$ws.Range("I92").Value = 933
$ws.Range("J92").Value = 2825.7144
$ws.Range("K92").Value = 2799
$ws.Range("L92").Value = 8477.143199999999
$ws.Range("M92").Value = -1551
$ws.Range("N92").Value = -10973.1432
$ws.Range("H132").Value = 3965
$ws.Range("I132").Value = 8399
$ws.Range("J132").Value = 2191.4
$ws.Range("K132").Value = 75591
$ws.Range("L132").Value = 19722.6
$ws.Range("M132").Value = -73061
$ws.Range("N132").Value = -24782.6
$ws.Range("H135").Value = 638
$ws.Range("I135").Value = 531.7143
$ws.Range("J135").Value = 720.6667
$ws.Range("K135").Value = 4785.428699999999
$ws.Range("L135").Value = 6486.0003
$ws.Range("M135").Value = -2250.428699999999
$ws.Range("N135").Value = -11556.0003

# Sheet: GSM (20 cells)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 50821.285
$ws.Range("J32").Value = 69249.75
$ws.Range("L32").Value = 69249.75
$ws.Range("N32").Value = -69841.75
$ws.Range("H42").Value = 73145
$ws.Range("J42").Value = 73145
$ws.Range("L42").Value = 73145
$ws.Range("N42").Value = -74115
$ws.Range("H102").Value = 9999.75
$ws.Range("I102").Value = 9999.75
$ws.Range("K102").Value = 9999.75
$ws.Range("M102").Value = -8377.75
$ws.Range("H115").Value = 73145
$ws.Range("J115").Value = 73145
$ws.Range("L115").Value = 73145
$ws.Range("N115").Value = -75495
$ws.Range("H132").Value = 5145.087
$ws.Range("I132").Value = 2225.7646
$ws.Range("K132").Value = 6677.293799999999
$ws.Range("M132").Value = -4147.293799999999

# Sheet: LTW (25 cells)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5378.222
$ws.Range("I46").Value = 1093.3334
$ws.Range("J46").Value = 5767.758
$ws.Range("K46").Value = 1093.3334
$ws.Range("L46").Value = 5767.758
$ws.Range("M46").Value = -905.3334
$ws.Range("N46").Value = -6143.758
$ws.Range("H55").Value = 1994.4
$ws.Range("I55").Value = 1998
$ws.Range("J55").Value = 1993.5
$ws.Range("K55").Value = 1998
$ws.Range("L55").Value = 1993.5
$ws.Range("M55").Value = -1825
$ws.Range("N55").Value = -2339.5
$ws.Range("H92").Value = 59000
$ws.Range("I92").Value = 40000
$ws.Range("J92").Value = 78000
$ws.Range("K92").Value = 40000
$ws.Range("L92").Value = 78000
$ws.Range("M92").Value = -37504
$ws.Range("N92").Value = -82992
$ws.Range("H93").Value = 2201
$ws.Range("I93").Value = 1933.3334
$ws.Range("K93").Value = 1933.3334
$ws.Range("M93").Value = -685.3334

# Sheet: WVR (19 cells)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 13999.223
$ws.Range("I58").Value = 10999.125
$ws.Range("J58").Value = 38000
$ws.Range("K58").Value = 10999.125
$ws.Range("L58").Value = 38000
$ws.Range("M58").Value = -10691.125
$ws.Range("N58").Value = -38616
$ws.Range("H122").Value = 4883.1665
$ws.Range("I122").Value = 6266.6665
$ws.Range("K122").Value = 18799.9995
$ws.Range("M122").Value = -16349.9995
$ws.Range("H132").Value = 10157
$ws.Range("I132").Value = 4999.5
$ws.Range("K132").Value = 14998.5
$ws.Range("M132").Value = -12468.5
$ws.Range("H136").Value = 4637.25
$ws.Range("I136").Value = 2580.875
$ws.Range("K136").Value = 7742.625
$ws.Range("M136").Value = -5192.625

Write-Host "Applied 213 cell updates across 8 sheets"